$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(45, 2).Value = 3

$ws.Cells.Item(46, 1).Value = "Ruilin"
$ws.Cells.Item(46, 2).Value = "'3"
$ws.Cells.Item(46, 3).Value = "无"
$ws.Cells.Item(46, 4).Value = "FBK"
$ws.Cells.Item(46, 5).Value = "WRI"
$ws.Cells.Item(46, 6).Value = "92b80f86-ee70-4a78-8469-1a9c33b052ed"
$ws.Cells.Item(46, 7).Value = "7Y52YHDS2X7ae_annotated.xlsx"
$ws.Cells.Item(46, 8).Value = "We will include this description in the new version of the paper."
